# Add a new worksheet "DTR-Features" at the end of the workbook, containing
# a pandas-style (index, feature, coefficient) table -- a Decision-Tree
# Regressor feature-importance ranking.

$wb = $excel.ActiveWorkbook

# Grab a reference cell that already carries the workbook's existing
# "bold + centered + thin border" header/index style (style index 1 in the
# original styles.xml) so we can clone it without minting new cellXfs
# entries.
$styleSource = $wb.Worksheets.Item("PCA_Top_Features").Range("A2")

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "DTR-Features"

# Move the new sheet to the very end of the tab order (after the current
# last sheet, which at this point is "PCA_Top_Features").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

# Re-resolve the worksheet object by name: after the reorder above, the
# original $newSheet reference can point at a stale index.
$ws = $wb.Worksheets.Item("DTR-Features")

# ---- Header row -----------------------------------------------------
$ws.Range("B1").Value = "Features"
$ws.Range("C1").Value = "Coefficients"

# ---- Data rows: index (A), feature name (B), coefficient (C) --------
# The coefficient column was exported from pandas as strings (e.g. "0.576"),
# so each value must land as text, not a number.
$data = @(
    @(5,  "degree_spondylolisthesis", "0.576"),
    @(3,  "sacral_slope",             "0.125"),
    @(4,  "pelvic_radius",            "0.103"),
    @(10, "sacrum_angle",             "0.071"),
    @(2,  "lumbar_lordosis_angle",    "0.063"),
    @(9,  "cervical_tilt",            "0.029"),
    @(1,  "pelvic_tilt",              "0.02"),
    @(8,  "thoracic_slope",           "0.011"),
    @(7,  "Direct_tilt",              "0.002"),
    @(0,  "pelvic_incidence",         "0.0"),
    @(6,  "pelvic_slope",             "0.0"),
    @(11, "scoliosis_slope",          "0.0")
)

# Scratch cell, well outside the used range, used to stage a text literal:
# writing ="0.576" as a formula and then pasting-as-values converts the
# formula's string result into a literal text cell (t="s") without Excel's
# usual "looks like a number -> store as number" auto-coercion, and without
# touching NumberFormat (which would otherwise mint a new style).
$scratch = $ws.Range("Z1")

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]

    $scratch.Formula = '="' + $item[2] + '"'
    $scratch.Copy()
    $ws.Cells.Item($row, 3).PasteSpecial(-4163)   # xlPasteValues

    $row++
}

$scratch.Clear()
$excel.CutCopyMode = $false

# ---- Styling ----------------------------------------------------------
# Clone the existing "bold, centered, thin-bordered" style (already used by
# every other sheet's header row + index column) onto B1:C1 and A2:A13,
# instead of building the formatting up property-by-property (which would
# mint brand-new style records).
$styleSource.Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A2:A13").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
